$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "87-77="
$t.Cell(1, 2).Range.Text = "65-57="
$t.Cell(1, 3).Range.Text = "7+78="
$t.Cell(1, 4).Range.Text = "30+64="
$t.Cell(1, 5).Range.Text = "62-46="

$t.Cell(2, 1).Range.Text = "75-63="
$t.Cell(2, 2).Range.Text = "99-15="
$t.Cell(2, 3).Range.Text = "76+18="
$t.Cell(2, 4).Range.Text = "44-14="
$t.Cell(2, 5).Range.Text = "52+11="

$t.Cell(3, 1).Range.Text = "11+16="
$t.Cell(3, 2).Range.Text = "32-28="
$t.Cell(3, 3).Range.Text = "78-47="
$t.Cell(3, 4).Range.Text = "86-6="
$t.Cell(3, 5).Range.Text = "82-69="

$t.Cell(4, 1).Range.Text = "6-1="
$t.Cell(4, 2).Range.Text = "54+35="
$t.Cell(4, 3).Range.Text = "5+63="
$t.Cell(4, 4).Range.Text = "61-51="
$t.Cell(4, 5).Range.Text = "19+37="

$t.Cell(5, 1).Range.Text = "55-50="
$t.Cell(5, 2).Range.Text = "96-26="
$t.Cell(5, 3).Range.Text = "79-78="
$t.Cell(5, 4).Range.Text = "20+17="
$t.Cell(5, 5).Range.Text = "37+30="

$t.Cell(6, 1).Range.Text = "65-10="
$t.Cell(6, 2).Range.Text = "50+34="
$t.Cell(6, 3).Range.Text = "43+1="
$t.Cell(6, 4).Range.Text = "32-24="
$t.Cell(6, 5).Range.Text = "70-15="

$t.Cell(7, 1).Range.Text = "34+39="
$t.Cell(7, 2).Range.Text = "86+7="
$t.Cell(7, 3).Range.Text = "59-10="
$t.Cell(7, 4).Range.Text = "23+26="
$t.Cell(7, 5).Range.Text = "28+46="

$t.Cell(8, 1).Range.Text = "69-68="
$t.Cell(8, 2).Range.Text = "93-57="
$t.Cell(8, 3).Range.Text = "5+74="
$t.Cell(8, 4).Range.Text = "32+28="
$t.Cell(8, 5).Range.Text = "79-20="

$t.Cell(9, 1).Range.Text = "56+41="
$t.Cell(9, 2).Range.Text = "37+14="
$t.Cell(9, 3).Range.Text = "82-49="
$t.Cell(9, 4).Range.Text = "98-84="
$t.Cell(9, 5).Range.Text = "2+30="

$t.Cell(10, 1).Range.Text = "85-32="
$t.Cell(10, 2).Range.Text = "17-11="
$t.Cell(10, 3).Range.Text = "60-59="
$t.Cell(10, 4).Range.Text = "74+12="
$t.Cell(10, 5).Range.Text = "85-5="

$t.Cell(11, 1).Range.Text = "40-4="
$t.Cell(11, 2).Range.Text = "1+45="
$t.Cell(11, 3).Range.Text = "39-24="
$t.Cell(11, 4).Range.Text = "17+8="
$t.Cell(11, 5).Range.Text = "13+33="

$t.Cell(12, 1).Range.Text = "51+38="
$t.Cell(12, 2).Range.Text = "37+31="
$t.Cell(12, 3).Range.Text = "36+8="
$t.Cell(12, 4).Range.Text = "36+60="
$t.Cell(12, 5).Range.Text = "2+55="

$t.Cell(13, 1).Range.Text = "83+16="
$t.Cell(13, 2).Range.Text = "26+0="
$t.Cell(13, 3).Range.Text = "66+0="
$t.Cell(13, 4).Range.Text = "96-46="
$t.Cell(13, 5).Range.Text = "99-39="

$t.Cell(14, 1).Range.Text = "61-21="
$t.Cell(14, 2).Range.Text = "53-21="
$t.Cell(14, 3).Range.Text = "94-31="
$t.Cell(14, 4).Range.Text = "28-14="
$t.Cell(14, 5).Range.Text = "45-8="

$t.Cell(15, 1).Range.Text = "79-14="
$t.Cell(15, 2).Range.Text = "41-21="
$t.Cell(15, 3).Range.Text = "55+15="
$t.Cell(15, 4).Range.Text = "2+15="
$t.Cell(15, 5).Range.Text = "63-7="

$t.Cell(16, 1).Range.Text = "59+0="
$t.Cell(16, 2).Range.Text = "97-10="
$t.Cell(16, 3).Range.Text = "53+2="
$t.Cell(16, 4).Range.Text = "70-38="
$t.Cell(16, 5).Range.Text = "97-23="

$t.Cell(17, 1).Range.Text = "12+79="
$t.Cell(17, 2).Range.Text = "17+23="
$t.Cell(17, 3).Range.Text = "9+62="
$t.Cell(17, 4).Range.Text = "95+1="
$t.Cell(17, 5).Range.Text = "15+2="

$t.Cell(18, 1).Range.Text = "24-7="
$t.Cell(18, 2).Range.Text = "14+49="
$t.Cell(18, 3).Range.Text = "58+1="
$t.Cell(18, 4).Range.Text = "10+44="
$t.Cell(18, 5).Range.Text = "4+72="

$t.Cell(19, 1).Range.Text = "90-5="
$t.Cell(19, 2).Range.Text = "79+17="
$t.Cell(19, 3).Range.Text = "45+51="
$t.Cell(19, 4).Range.Text = "20+75="
$t.Cell(19, 5).Range.Text = "97-42="

$t.Cell(20, 1).Range.Text = "38-17="
$t.Cell(20, 2).Range.Text = "71-23="
$t.Cell(20, 3).Range.Text = "66+21="
$t.Cell(20, 4).Range.Text = "7+6="
$t.Cell(20, 5).Range.Text = "62+33="
